# The deck's slide master ("Design 1") currently uses the custom
# "Integral" theme (ppt/theme/theme1.xml). The commit swaps the
# presentation's design over to the plain built-in "Office Theme"
# palette (the one previously only used by the Notes Master).
#
# The PowerPoint object model only exposes the *applied* (slide-master)
# theme for editing -- via ThemeColorScheme.Colors(i).RGB -- so we
# recolor it to the Office Theme's 12 theme colors. RGB() longs are
# R + G*256 + B*65536 (PowerPoint's native color-value packing).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Theme color order (1-based, matches PowerPoint's ThemeColorScheme.Colors):
#  1 Dark1 / Text1        2 Light1 / Background1
#  3 Dark2 / Text2        4 Light2 / Background2
#  5 Accent1               6 Accent2
#  7 Accent3               8 Accent4
#  9 Accent5               10 Accent6
# 11 Hyperlink             12 FollowedHyperlink

$colors.Item(1).RGB  = 0         # Dark1 / Text1        -> 000000
$colors.Item(2).RGB  = 16777215  # Light1 / Background1 -> FFFFFF
$colors.Item(3).RGB  = 6968388   # Dark2 / Text2         -> 44546A
$colors.Item(4).RGB  = 15132391  # Light2 / Background2  -> E7E6E6
$colors.Item(5).RGB  = 13998939  # Accent1               -> 5B9BD5
$colors.Item(6).RGB  = 3243501   # Accent2               -> ED7D31
$colors.Item(7).RGB  = 10855845  # Accent3               -> A5A5A5
$colors.Item(8).RGB  = 49407     # Accent4               -> FFC000
$colors.Item(9).RGB  = 12874308  # Accent5               -> 4472C4
$colors.Item(10).RGB = 4697456   # Accent6               -> 70AD47
$colors.Item(11).RGB = 12673797  # Hyperlink             -> 0563C1
$colors.Item(12).RGB = 7491477   # FollowedHyperlink     -> 954F72
